$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new data block: A3 = B3+C3 (formula), B3 = 1, C3 = 2
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 2
$ws.Range("A3").Formula = "=B3+C3"

# New conditional format rule on A3: highlight when value > 10
# (1 = xlCellValue, 5 = xlGreater)
$rule = $ws.Range("A3").FormatConditions.Add(1, 5, "10")
$rule.Font.Color = 393372       # dark red FF9C0006, BGR-packed for COM
$rule.Interior.Color = 13551615 # light red/pink fill FFFFC7CE, BGR-packed
$rule.SetFirstPriority()        # new rule evaluates before the A1:B1 rule

# Select A3, matching the new active selection in the saved file
$ws.Range("A3").Select() | Out-Null
